$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column C ("Förändrad") date serials from 45185 to 45204 for rows 2-13
$ws.Range("C2:C13").Value = 45204
